$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared string used as a note for the added row ---
$note = "Add Memory Control 2 and ROM 2"

# --- Row 45: new timesheet entry (VGA Top / Arch) ---

# Column A holds the date as literal text (not a real date serial). A
# direct string assignment would get auto-converted to a date serial by
# Excel, so instead enter it as a text-producing formula (never
# reinterpreted as a date) and then flatten it down to its literal
# result, which keeps it a plain text value sharing the existing
# "4.4.2020" string already used elsewhere in the sheet.
$ws.Range("A45").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight
$ws.Range("A45").VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$ws.Range("A45").Formula = '="4.4.2020"'
$ws.Range("A45").Copy()
$ws.Range("A45").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

# Columns B/C hold the From/To times, formatted as h:mm.
$ws.Range("B45:C45").NumberFormat = "h:mm"
$ws.Range("B45").Value = 0.5
$ws.Range("C45").Value = 0.51041666666666663

# Column D computes the elapsed duration, also formatted as h:mm.
$ws.Range("D45").NumberFormat = "h:mm"
$ws.Range("D45").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight
$ws.Range("D45").VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$ws.Range("D45").Formula = "=C45-B45"

# Columns E/F hold the Task/Unit text.
$ws.Range("E45:F45").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft
$ws.Range("E45:F45").VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$ws.Range("E45").Value = "VGA Top"
$ws.Range("F45").Value = "Arch"

# Column G holds the free-form note describing the change.
$ws.Range("G45").Value = $note

# Move the active selection to A46, just past the newly added row.
$ws.Range("A46").Select()
